$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update NOMBRE value
$ws.Range("B2").Value = "Cerveza Artesanal Actualizada test"

# Update DESCRIPCIÓN value
$ws.Range("B3").Value = "testttasdasd"

# Update MARCA value
$ws.Range("B4").Value = "carozzitest"

# Update CONTENIDO value - "400" looks numeric, so a plain .Value assignment
# would store it as a number. Round-trip it through a formula + copy/paste
# special so it lands back in the cell as text (matches original text-typed
# shared-string cell and keeps the existing cell style).
$helper = $ws.Cells.Item(20, 5)
$helper.Formula = "=""400"""
$helper.Copy()
$ws.Range("B5").PasteSpecial(-4163)
$helper.Clear()

# Update PRECIO value
$ws.Range("B8").Value = "NO REGISTRADO"

# Update CATEGORÍA value
$ws.Range("B9").Value = "Vino"
